$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 28 (row 40) results
$ws.Range("E40").Value = 0
$ws.Range("H40").Value = 60
$ws.Range("K40").Value = 70
$ws.Range("N40").Value = 100
$ws.Range("Q40").Value = 20
$ws.Range("T40").Value = 40
$ws.Range("W40").Value = 50
$ws.Range("Z40").Value = 80
$ws.Range("AC40").Value = 30

# Contest 29 (row 41) results
$ws.Range("E41").Value = 80
$ws.Range("H41").Value = 0
$ws.Range("K41").Value = 50
$ws.Range("N41").Value = 40
$ws.Range("Q41").Value = 20
$ws.Range("T41").Value = 70
$ws.Range("W41").Value = 100
$ws.Range("Z41").Value = 30
$ws.Range("AC41").Value = 60
